# adj bpsk 1200 tune
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 33: update H33 formula to reference the second ($B$29/$B$30) block ---
$ws.Range("H33").Formula = "=ROUND(((F33*(1+G33)) * `$B`$29/`$B`$30)/B33, 0)"

# --- Fill in the plain (non-formula) inputs for rows 34-36 first ---
$ws.Range("F34").Value = 25
$ws.Range("G34").Value = 0.5
$ws.Range("I34").Value = 0.81499999999999995

$ws.Range("B35").Value = 0.00021000000000000001
$ws.Range("D35").Value = 0.2
$ws.Range("F35").Value = 50
$ws.Range("G35").Value = 0.5
$ws.Range("I35").Value = 0.81499999999999995

$ws.Range("F36").Value = 25
$ws.Range("G36").Value = 0.5
$ws.Range("I36").Value = 0.81499999999999995

# --- Assign the formulas across the whole 34:36 ranges at once so Excel
#     stores them as shared formulas (matching the original author's fill-down) ---
$ws.Range("C34:C36").Formula = "=ROUND(POWER(2,`$C`$22)*B34, 0)"
$ws.Range("E34:E36").Formula = "=ROUND(POWER(2,`$E`$22)*D34, 0)"
$ws.Range("H34:H36").Formula = "=ROUND(((F34*(1+G34)) * `$B`$29/`$B`$30)/B34, 0)"
$ws.Range("J34:J36").Formula = "=LOG(H34*C34,2)"

$wb.Application.Calculate()

# --- Update the saved selection to H35 ---
$ws.Range("H35").Select()
